$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the steel ("S") description in the industrial mapping string:
# remove the "/RME" segment from the "28% S/LFM+CDL/RME/H:1" line.
$cell = $ws.Range("B2")
$text = $cell.Value()
$text = $text.Replace("28% S/LFM+CDL/RME/H:1", "28% S/LFM+CDL/H:1")

# Wrap the long multi-line description and update the cell value.
$cell.WrapText = $true
$cell.Value = $text

# Expand the row so the wrapped text is fully visible.
$ws.Rows.Item(2).RowHeight = 409.6

# Restore the selection left behind by the edit.
$ws.Range("B13").Select() | Out-Null
